# ---------------------------------------------------------------------------
# SmartDietAPI dishes.xlsx - "config mapping with crud in frontend"
#
# Adds three new dish rows (plus two blank "x" placeholder rows copied from
# the template row) to the "Meal_Data" table, wires up their hyperlinks,
# and grows every range-bound artifact (table, autofilter, filter defined
# name, data validations, frozen-pane selection) from a single-row template
# (A1:K2) out to the new seven-row extent (A1:K7).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Cell values for the new rows.
# ---------------------------------------------------------------------------

# Row 2 - "Ba roi rim chay canh"
$ws.Range("A2").Value = "Ba rọi rim cháy cạnh"
$ws.Range("B2").Value = "Nam"
$ws.Range("C2").Value = "Giàu protein"
$d2 = @'
Ba rọi rim cháy cạnh ngon bắt cơm với nước rim khô kẹo, miếng thịt cháy cạnh rám mặt mềm bên trong.
Có màu caramel đặc trưng của các món kho
Dậy mùi kho đặc trưng của món Việt
'@
$ws.Range("D2").Value = $d2
$ws.Range("H2").Value = "Chờ kiểm tra"
$ws.Range("J2").Value = "Huỳnh Nhật"

# Row 3 - "Bach Tuoc Ngam Sa Tac Chua Cay"
$ws.Range("A3").Value = "Bạch Tuộc Ngâm Sả Tắc Chua Cay"
$ws.Range("B3").Value = "Chung"
$ws.Range("C3").Value = "Giàu calo"
$ws.Range("D3").Value = " Từng con bạch tuộc sau khi hấp xong mềm, dai dai, không có vị tanh nhờ có gừng. Sả, tắc, ớt kết hợp làm cho món ăn trở nên đẹp mắt, thu hút."
$ws.Range("H3").Value = "Chờ kiểm tra"
$ws.Range("J3").Value = "Linh Chi"

# Row 4 - "Banh bao sup - Xiao long bao"
$ws.Range("A4").Value = "Bánh bao súp - Xiao long bao"
$ws.Range("B4").Value = "Chung"
$ws.Range("C4").Value = "Giàu calo"
$ws.Range("D4").Value = "Cách làm bánh bao hấp đơn giản với lớp vỏ bánh dai mỏng bên ngoài, bao bọc nhân thịt, nước súp chảy ra từ bên trong bánh rất đặc trưng và hấp dẫn. "
$ws.Range("H4").Value = "Chờ kiểm tra"
$ws.Range("J4").Value = "Linh Chi"

# Row 5 - duplicate of row 4 (same dish entered twice)
$ws.Range("A5").Value = "Bánh bao súp - Xiao long bao"
$ws.Range("B5").Value = "Chung"
$ws.Range("C5").Value = "Giàu calo"
$ws.Range("D5").Value = "Cách làm bánh bao hấp đơn giản với lớp vỏ bánh dai mỏng bên ngoài, bao bọc nhân thịt, nước súp chảy ra từ bên trong bánh rất đặc trưng và hấp dẫn. "
$ws.Range("H5").Value = "Chờ kiểm tra"
$ws.Range("J5").Value = "Linh Chi"

# Rows 6 & 7 - blank placeholder rows (template row filled down, marked "x")
$ws.Range("A6").Value = "x"
$ws.Range("A7").Value = "x"

Write-Host "values set"

# ---------------------------------------------------------------------------
# 2. Hyperlinks (Link công thức / Nguồn tham khảo / Link ảnh columns).
# ---------------------------------------------------------------------------

$ws.Hyperlinks.Add($ws.Range("F2"), "https://docs.google.com/document/d/ba-roi-rim-chay-canh", "", "", "Link công thức - Google Tài liệu") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G2"), "https://monngonmoingay.com/ba-roi-rim-chay-canh", "", "", "Ba rọi rim cháy cạnh | Món Ngon Mỗi Ngày") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I2"), "https://drive.google.com/drive/folders/M134-M204", "", "", "M134-M204 - Google Drive") | Out-Null

$ws.Hyperlinks.Add($ws.Range("F3"), "https://docs.google.com/document/d/chi-thong-tin-mon-an-bach-tuoc", "", "", "Chi_thông tin món ăn") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G3"), "https://www.cooky.vn/cong-thuc/bach-tuoc-ngam-sa-tac-chua-cay-18935", "", "", "https://www.cooky.vn/cong-thuc/bach-tuoc-ngam-sa-tac-chua-cay-18935") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I3"), "https://drive.google.com/file/d/M104-jpeg", "", "", "M104.jpeg") | Out-Null

$ws.Hyperlinks.Add($ws.Range("F4"), "https://docs.google.com/document/d/chi-thong-tin-mon-an-banh-bao-sup", "", "", "Chi_thông tin món ăn") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G4"), "https://www.cooky.vn/cong-thuc/banh-bao-sup-xiao-long-bao-34747", "", "", "https://www.cooky.vn/cong-thuc/banh-bao-sup-xiao-long-bao-34747") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I4"), "https://drive.google.com/file/d/M113-jpeg", "", "", "M113.jpeg") | Out-Null

$ws.Hyperlinks.Add($ws.Range("F5"), "https://docs.google.com/document/d/chi-thong-tin-mon-an-banh-bao-sup", "", "", "Chi_thông tin món ăn") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G5"), "https://www.cooky.vn/cong-thuc/banh-bao-sup-xiao-long-bao-34747", "", "", "https://www.cooky.vn/cong-thuc/banh-bao-sup-xiao-long-bao-34747") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I5"), "https://drive.google.com/file/d/M113-jpeg", "", "", "M113.jpeg") | Out-Null

Write-Host "hyperlinks set"

# ---------------------------------------------------------------------------
# 3. Grow the table / autofilter from A1:K2 to A1:K7.
# ---------------------------------------------------------------------------

$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:K7")) | Out-Null

Write-Host "table resized"

# ---------------------------------------------------------------------------
# 4. Grow the _FilterDatabase defined name to match.
# ---------------------------------------------------------------------------

$fdb = $wb.Names.Item(1)
$fdb.RefersTo = "='Món ăn'!`$A`$1:`$K`$7"

Write-Host "defined name updated"

# ---------------------------------------------------------------------------
# 5. Grow every data-validation sqref from row 2 to rows 2:7.
# ---------------------------------------------------------------------------

$ws.Range("J2").Validation.Delete()
$rJ = $ws.Range("J2:J7")
$rJ.Validation.Add(3, 1, 1, '"Huỳnh Nhật,Linh Chi,Anh Khoa"') | Out-Null
$rJ.Validation.IgnoreBlank = $true
$rJ.Validation.InCellDropdown = $true
$rJ.Validation.ShowInput = $false
$rJ.Validation.ShowError = $false

$ws.Range("E2").Validation.Delete()
$rE = $ws.Range("E2:E7")
$rE.Validation.Add(3, 1, 1, '"Option 1,Option 2"') | Out-Null
$rE.Validation.IgnoreBlank = $true
$rE.Validation.InCellDropdown = $true
$rE.Validation.ShowInput = $false
$rE.Validation.ShowError = $false

$ws.Range("B2").Validation.Delete()
$rB = $ws.Range("B2:B7")
$rB.Validation.Add(3, 1, 1, '"Bắc,Trung,Nam,Chung"') | Out-Null
$rB.Validation.IgnoreBlank = $true
$rB.Validation.InCellDropdown = $true
$rB.Validation.ShowInput = $false
$rB.Validation.ShowError = $false

$ws.Range("H2").Validation.Delete()
$rH = $ws.Range("H2:H7")
$rH.Validation.Add(3, 1, 1, '"Chờ kiểm tra,Đợi Feedback,Đã Feedback,Đã xong"') | Out-Null
$rH.Validation.IgnoreBlank = $true
$rH.Validation.InCellDropdown = $true
$rH.Validation.ShowInput = $false
$rH.Validation.ShowError = $false

$ws.Range("C2").Validation.Delete()
$rC = $ws.Range("C2:C7")
$rC.Validation.Add(3, 1, 1, '"Không,Lỏng,Ít calo,Giàu calo,Ít cholesterol,Ăn chay,Ít natri,Giàu protein,Ít protein"') | Out-Null
$rC.Validation.IgnoreBlank = $true
$rC.Validation.InCellDropdown = $true
$rC.Validation.ShowInput = $false
$rC.Validation.ShowError = $false

# Plain "allow blank / show dropdown button" validation (no list) - originally
# sqref="A2 K2 I2 F2:G2 D2", now A2:A7 D2:D7 F2:G7 I2:I7 K2:K7.
$ws.Range("A2").Validation.Delete()
$ws.Range("K2").Validation.Delete()
$ws.Range("I2").Validation.Delete()
$ws.Range("F2:G2").Validation.Delete()
$ws.Range("D2").Validation.Delete()
foreach ($addr in @("A2:A7", "D2:D7", "F2:G7", "I2:I7", "K2:K7")) {
    $r = $ws.Range($addr)
    $r.Validation.Add(0) | Out-Null
    $r.Validation.IgnoreBlank = $true
    $r.Validation.InCellDropdown = $true
    $r.Validation.ShowInput = $false
    $r.Validation.ShowError = $false
}

Write-Host "data validations grown"
